# Updates cryptocurrency price/volume figures (and re-orders the NEAR
# Protocol / LidoDAOToken rows) on Sheet1, per the Thu Dec 28 2023 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.746.50"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "2.372.48"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  -0.07%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'333.29"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +8.26%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").Value = "'100.66"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -5.67%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  +0.10%  "
$style = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.637"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  +1.38%  "
$style = $ws.Range("D10").Style
$ws.Range("D10").Value = "'40.22"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -4.54%  "
$ws.Range("E11").Value = "  -1.23%  "
$style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'8.51"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -3.94%  "
$ws.Range("E13").Value = "  -2.88%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'0.106"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  +0.52%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'16.36"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "2.729.19"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").Value = "2.371.66"
$ws.Range("E17").Value = "  -3.29%  "
$ws.Range("D18").Value = "42.673.66"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("E19").Value = "  +8.36%  "
$ws.Range("E20").Value = "  -0.71%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").Value = "'3.83"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +13.05%  "
$style = $ws.Range("D22").Style
$ws.Range("D22").Value = "'75.58"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  +0.58%  "
$style = $ws.Range("D23").Style
$ws.Range("D23").Value = "'272.37"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  +8.01%  "
$ws.Range("E24").Value = "  -6.72%  "
$style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'9.85"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +10.74%  "
$ws.Range("E26").Value = "  -0.10%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").Value = "'11.51"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -3.48%  "
$ws.Range("E28").Value = "  -0.22%  "
$style = $ws.Range("D29").Style
$ws.Range("D29").Value = "'23.38"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  +3.62%  "
$ws.Range("E30").Value = "  -1.33%  "
$style = $ws.Range("D31").Style
$ws.Range("D31").Value = "'175.30"
$ws.Range("D31").Style = $style
$style = $ws.Range("D32").Style
$ws.Range("D32").Value = "'3.09"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -2.65%  "
$style = $ws.Range("D33").Style
$ws.Range("D33").Value = "'0.0910"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -0.86%  "
$style = $ws.Range("D34").Style
$ws.Range("D34").Value = "'35.53"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -8.67%  "
$style = $ws.Range("D35").Style
$ws.Range("D35").Value = "'6.12"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  +2.83%  "
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("E37").Value = "  -8.68%  "
$ws.Range("E38").Value = "  -3.99%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$style = $ws.Range("D39").Style
$ws.Range("D39").Value = "'3.88"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -4.14%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$style = $ws.Range("D40").Style
$ws.Range("D40").Value = "'2.91"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  +6.63%  "
$ws.Range("E41").Value = "  +3.54%  "
$style = $ws.Range("D42").Style
$ws.Range("D42").Value = "'1.53"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  +2.21%  "
$style = $ws.Range("D43").Style
$ws.Range("D43").Value = "'0.235"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  +2.65%  "
$style = $ws.Range("D44").Style
$ws.Range("D44").Value = "'70.37"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -1.82%  "
$ws.Range("E45").Value = "  +0.11%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'118.38"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +7.50%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'12.13"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -1.54%  "
$style = $ws.Range("D48").Style
$ws.Range("D48").Value = "'89.12"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +28.93%  "
$style = $ws.Range("D49").Style
$ws.Range("D49").Value = "'5.50"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -1.87%  "
$style = $ws.Range("D50").Style
$ws.Range("D50").Value = "'9.15"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  -1.89%  "
$ws.Range("D51").Value = "1.584.91"
$ws.Range("E51").Value = "  +6.07%  "
